# Update cryptocurrency price (column D) and volume-change (column E) figures
# to reflect the latest scrape, per the GitHub Actions "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.256.48"
$ws.Range("E2").Value = "  +2.52%  "

# Row 3
$ws.Range("D3").Value = "3.566.70"
$ws.Range("E3").Value = "  +5.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.75"
$ws.Range("E5").Value = "  +0.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.83"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
$ws.Range("D7").Value = "3.556.81"
$ws.Range("E7").Value = "  +4.99%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +1.66%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.206"
$ws.Range("E10").Value = "  +5.68%  "

# Row 11
$ws.Range("E11").Value = "  +2.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.93"
$ws.Range("E12").Value = "  +2.67%  "

# Row 13
$ws.Range("E13").Value = "  +2.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "696.29"
$ws.Range("E14").Value = "  +2.16%  "

# Row 15
$ws.Range("D15").Value = "4.139.24"
$ws.Range("E15").Value = "  +5.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.87"
$ws.Range("E16").Value = "  +2.56%  "

# Row 17
$ws.Range("D17").Value = "71.381.04"
$ws.Range("E17").Value = "  +2.70%  "

# Row 18
$ws.Range("D18").Value = "3.534.79"
$ws.Range("E18").Value = "  +4.11%  "

# Row 19
$ws.Range("E19").Value = "  +1.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.23"
$ws.Range("E20").Value = "  +3.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.66"
$ws.Range("E21").Value = "  +3.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.928"
$ws.Range("E22").Value = "  +2.56%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.51"
$ws.Range("E23").Value = "  +1.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.55"
$ws.Range("E24").Value = "  +2.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "104.02"
$ws.Range("E25").Value = "  +0.29%  "

# Row 26
$ws.Range("E26").Value = "  +1.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +1.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  +2.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.36"
$ws.Range("E29").Value = "  +3.46%  "

# Row 30
$ws.Range("E30").Value = "  +3.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.44"
$ws.Range("E31").Value = "  +5.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.06"
$ws.Range("E32").Value = "  +12.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "587.92"
$ws.Range("E33").Value = "  +4.99%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.26"
$ws.Range("E34").Value = "  +1.02%  "

# Row 35
$ws.Range("E35").Value = "  -0.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.56"
$ws.Range("E36").Value = "  +2.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"

# Row 38
$ws.Range("D38").Value = "3.660.14"
$ws.Range("E38").Value = "  -0.91%  "

# Row 39
$ws.Range("E39").Value = "  +3.65%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.95"
$ws.Range("E40").Value = "  +0.96%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0762"
$ws.Range("E41").Value = "  +9.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.45"
$ws.Range("E42").Value = "  +5.15%  "

# Row 43
$ws.Range("E43").Value = "  +2.68%  "

# Row 44
$ws.Range("E44").Value = "  +3.30%  "

# Row 45
$ws.Range("E45").Value = "  +1.52%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +2.09%  "

# Row 47
$ws.Range("E47").Value = "  +2.31%  "

# Row 48
$ws.Range("E48").Value = "  +5.13%  "

# Row 49
$ws.Range("E49").Value = "  +1.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.23"
$ws.Range("E51").Value = "  +0.87%  "
